$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FullSuite")

# Copy formatting (fill + wrap + alignment) from row 245 down to the new rows 246:256
$ws.Range("A245:C245").Copy()
$ws.Range("A246:C256").PasteSpecial(-4122)

# Row 246: description + result text, set together (matches shared-string insertion order)
$ws.Cells.Item(246,1).Value = "Payroll Suite ScottishTaxMonth1CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(246,2).Value = "executes ScottishTaxMonth1CSBRNTK50PercentRegulatory201819 scenario"

# Column A (Description) for rows 247-256
$ws.Cells.Item(247,1).Value = "Payroll Suite ScottishTaxMonth2CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(248,1).Value = "Payroll Suite ScottishTaxMonth3CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(249,1).Value = "Payroll Suite ScottishTaxMonth4CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(250,1).Value = "Payroll Suite ScottishTaxMonth5CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(251,1).Value = "Payroll Suite ScottishTaxMonth6CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(252,1).Value = "Payroll Suite ScottishTaxMonth7CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(253,1).Value = "Payroll Suite ScottishTaxMonth8CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(254,1).Value = "Payroll Suite ScottishTaxMonth9CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(255,1).Value = "Payroll Suite ScottishTaxMonth10CSBRNTK50PercentRegulatory201819"
$ws.Cells.Item(256,1).Value = "Payroll Suite ScottishTaxMonth11CSBRNTK50PercentRegulatory201819"

# Column B (Results) for rows 247-256
$ws.Cells.Item(247,2).Value = "executes ScottishTaxMonth2CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(248,2).Value = "executes ScottishTaxMonth3CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(249,2).Value = "executes ScottishTaxMonth4CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(250,2).Value = "executes ScottishTaxMonth5CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(251,2).Value = "executes ScottishTaxMonth6CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(252,2).Value = "executes ScottishTaxMonth7CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(253,2).Value = "executes ScottishTaxMonth8CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(254,2).Value = "executes ScottishTaxMonth9CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(255,2).Value = "executes ScottishTaxMonth10CSBRNTK50PercentRegulatory201819 scenario"
$ws.Cells.Item(256,2).Value = "executes ScottishTaxMonth11CSBRNTK50PercentRegulatory201819 scenario"

# Column C (Runmode) for rows 246-256
For ($r = 246; $r -le 256; $r++) {
    $ws.Cells.Item($r,3).Value = "Y"
}

# Row heights
$ws.Rows.Item(244).RowHeight = 30.75
$ws.Rows.Item(245).RowHeight = 33.75
$ws.Rows.Item(246).RowHeight = 45.75
For ($r = 247; $r -le 256; $r++) {
    $ws.Rows.Item($r).RowHeight = 45
}

# Column B width
$ws.Columns.Item(2).ColumnWidth = 37.67

# View / selection
$ws.Range("G253").Select()
